$d = $word.ActiveDocument

# 1. "Front-end Development , Bootstrap" -> "Front-end Development, Microsoft Office"
#    a) collapse "Development , " into "Development, " (removes the stray space before the comma)
$d.Content.Find.Execute("Development ,", $true, $false, $false, $false, $false, $true, 1, $false, "Development,", 2)

#    b) swap the certification name
$d.Content.Find.Execute("Bootstrap", $true, $false, $false, $false, $false, $true, 1, $false, "Microsoft Office", 2)

# 2. Update the signature date from 19-05-2023 to 03-07-2023
$d.Content.Find.Execute("19-05-2023", $true, $false, $false, $false, $false, $true, 1, $false, "03-07-2023", 2)

$d.Save()
